# Apply updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "60.654.74"
$ws.Cells.Item(2, 5).Value = "  -0.45%  "
$ws.Cells.Item(3, 4).Value = "2.401.29"
$ws.Cells.Item(3, 5).Value = "  -0.80%  "
$ws.Cells.Item(4, 5).Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.15"
$ws.Range("D6").Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.71%  "
$ws.Cells.Item(7, 5).Value = "  -0.34%  "
$ws.Cells.Item(8, 5).Value = "  +1.42%  "
$ws.Cells.Item(9, 4).Value = "2.408.20"
$ws.Cells.Item(9, 5).Value = "  +0.13%  "
$ws.Cells.Item(10, 5).Value = "  -0.26%  "
$ws.Cells.Item(11, 5).Value = "  -0.85%  "
$ws.Cells.Item(12, 5).Value = "  +1.67%  "
$ws.Cells.Item(13, 5).Value = "  +0.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.08"
$ws.Range("D14").Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000168"
$ws.Range("D15").Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -2.01%  "
$ws.Cells.Item(16, 4).Value = "2.813.05"
$ws.Cells.Item(16, 5).Value = "  -0.51%  "
$ws.Cells.Item(17, 4).Value = "60.593.24"
$ws.Cells.Item(17, 5).Value = "  -0.38%  "
$ws.Cells.Item(18, 4).Value = "2.405.09"
$ws.Cells.Item(18, 5).Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.10"
$ws.Range("D19").Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +9.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.63"
$ws.Range("D20").Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.19"
$ws.Range("D21").Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.23%  "
$ws.Cells.Item(22, 5).Value = "  +0.55%  "
$ws.Cells.Item(23, 5).Value = "  +0.11%  "
$ws.Cells.Item(24, 5).Value = "  -0.13%  "
$ws.Cells.Item(25, 5).Value = "  -2.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.09"
$ws.Range("D26").Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "569.08"
$ws.Range("D27").Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -1.68%  "
$ws.Cells.Item(28, 5).Value = "  -5.39%  "
$ws.Cells.Item(30, 4).Value = "0.0₃0934"
$ws.Cells.Item(30, 5).Value = "  +1.35%  "
$ws.Cells.Item(31, 5).Value = "  +2.02%  "
$ws.Cells.Item(32, 5).Value = "  -1.14%  "
$ws.Cells.Item(33, 5).Value = "  -1.37%  "
$ws.Cells.Item(34, 5).Value = "  -1.77%  "
$ws.Cells.Item(35, 5).Value = "  -0.61%  "
$ws.Cells.Item(36, 5).Value = "  +4.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.27"
$ws.Range("D37").Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +2.42%  "
$ws.Cells.Item(38, 5).Value = "  +0.62%  "
$ws.Cells.Item(39, 5).Value = "  -1.81%  "
$ws.Cells.Item(40, 5).Value = "  +0.11%  "
$ws.Cells.Item(41, 5).Value = "  -0.45%  "
$ws.Cells.Item(42, 5).Value = "  -0.09%  "
$ws.Cells.Item(43, 2).Value = "OKB"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.68"
$ws.Range("D43").Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +1.32%  "
$ws.Cells.Item(44, 5).Value = "  +0.46%  "
$ws.Cells.Item(45, 2).Value = "dogwifhat"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.49"
$ws.Range("D45").Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +5.92%  "
$ws.Cells.Item(46, 5).Value = "  +1.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.22"
$ws.Range("D47").Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.55"
$ws.Range("D48").Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.591"
$ws.Range("D49").Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0506"
$ws.Range("D50").Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.34"
$ws.Range("D51").Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.31%  "
